$wb = $excel.ActiveWorkbook

# Rows (1-based, matching worksheet row numbers) whose Priority column should
# be set to "ht" on both the zh-cn and de-de handback-status sheets.
$rows = @(7, 9, 10, 11, 12, 13)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $ws.Range("E$r").Value = "ht"
    }
}

# Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps to reflect the newly generated handoff report.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-12 20:24:09"
}

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-08-12 20:23:57"
}

# de-de's handoff run finished at the same moment captured by the Overview
# sheet's "Latest HO Xliff Generate Date" column.
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-08-12 20:24:09"
}
